$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was M, now B)
$ws.Range("A2").Value = "B"
$ws.Range("B2").Value = 0.9111111111111111
$ws.Range("C2").Value = 0.9534883720930233
$ws.Range("D2").Value = 0.9318181818181818
$ws.Range("E2").Value = 215

# Row 3 (was B, now M)
$ws.Range("A3").Value = "M"
$ws.Range("B3").Value = 0.9145299145299145
$ws.Range("C3").Value = 0.84251968503937
$ws.Range("D3").Value = 0.8770491803278688
$ws.Range("E3").Value = 127

# Row 5 (macro avg) - updated values
$ws.Range("B5").Value = 0.9128205128205128
$ws.Range("C5").Value = 0.8980040285661967
$ws.Range("D5").Value = 0.9044336810730254

# Row 6 (weighted avg) - updated values
$ws.Range("B6").Value = 0.9123806667666318
$ws.Range("D6").Value = 0.9114799853583289
